# Rename the inline picture objects that live in the document's headers
# and footers:
#   - the BTEC logo pictures (in the two headers) go from "image1.jpg"
#     to "image2.jpg"
#   - the Pearson logo pictures (in the two footers) go from
#     "image2.png" to "image1.png"
#
# Word's InlineShape object has no documented/settable .Name in the real
# object model either; renaming an inline picture's docPr/name normally
# requires going Range.Select() -> Selection.InlineShapes(1).Name = ...
# (some HeaderFooter ranges otherwise report a stale/unaddressed block).

$d = $word.ActiveDocument

function Set-InlinePictureName($range, $newName) {
    $range.Select()
    $sel = $word.Selection
    if ($sel.InlineShapes.Count -gt 0) {
        $sel.InlineShapes.Item(1).Name = $newName
    }
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    # Headers: BTEC logo, image1.jpg -> image2.jpg
    for ($h = 1; $h -le 3; $h++) {
        $hdr = $sec.Headers.Item($h)
        if ($hdr.Exists) {
            $hdrRange = $hdr.Range
            if ($hdrRange.InlineShapes.Count -gt 0) {
                $shp = $hdrRange.InlineShapes.Item(1)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    Set-InlinePictureName $shp.Range "image2.jpg"
                }
            }
        }
    }

    # Footers: Pearson logo, image2.png -> image1.png
    for ($f = 1; $f -le 3; $f++) {
        $ftr = $sec.Footers.Item($f)
        if ($ftr.Exists) {
            $ftrRange = $ftr.Range
            if ($ftrRange.InlineShapes.Count -gt 0) {
                $shp = $ftrRange.InlineShapes.Item(1)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    Set-InlinePictureName $shp.Range "image1.png"
                }
            }
        }
    }
}
